$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (column D) values
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "264.71"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "22.73"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "6.270"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06150"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.595"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.713"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.349"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8297"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.01356"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1581"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08215"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03403"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03144"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.915"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.001712"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.04795"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006258"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.005932"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.001102"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.771"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.306"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3378"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04636"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006954"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1136"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01034"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006177"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.7782"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.2020"

# Update Hora (column G) values: all rows 2-51 go from 6 to 7
for ($row = 2; $row -le 51; $row++) {
    $cell = $ws.Range("G" + $row)
    $cell.NumberFormat = "@"
    $cell.Value = "7"
}

